# Apply the cell value changes described by the diff to the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value  = "-"
$ws.Range("C3").Value  = "[-, 'MEC-3B-Ens. Dest. Não Desti.', -, -]"
$ws.Range("C4").Value  = "[-, 'MEC-3B-Ens. Dest. Não Desti.', -, -]"
$ws.Range("C6").Value  = "[-, 'MEC-3B-Ens. Dest. Não Desti.', -, -]"
$ws.Range("E6").Value  = "-"
$ws.Range("E7").Value  = "-"
$ws.Range("E10").Value = "-"
$ws.Range("F11").Value = "['MEC-3A-Ens. Dest. Não Desti.', -, -, -]"
$ws.Range("F12").Value = "['MEC-3A-Ens. Dest. Não Desti.', -, -, -]"
$ws.Range("E16").Value = "-"
$ws.Range("C18").Value = "['MEC-1NA-E. D. N. D.', -, -, -]"
$ws.Range("D18").Value = "-"
$ws.Range("F18").Value = "[-, -, 'MEC-1NB-E. D. N. D.', -]"
$ws.Range("C19").Value = "['MEC-1NA-E. D. N. D.', -, -, -]"
$ws.Range("C20").Value = "['MEC-1NA-E. D. N. D.', -, -, -]"
$ws.Range("D20").Value = "-"
$ws.Range("F20").Value = "[-, -, 'MEC-1NB-E. D. N. D.', -]"
$ws.Range("C21").Value = "['MEC-1NA-E. D. N. D.', -, -, -]"
$ws.Range("D21").Value = "-"
$ws.Range("F21").Value = "[-, -, 'MEC-1NB-E. D. N. D.', -]"
